# Update countries & provincias Spain
# - Refresh "datos actualizados" timestamp (row 1)
# - Swap Barein/Bolivia data (rows 50-51) with refreshed counts
# - Swap Marruecos/Honduras data (rows 67-68) with refreshed counts
# - Refresh Haiti counts (row 82)
# - Refresh Nueva Caledonia counts (row 199)
# - Swap Santa Sede/Islas Turcas y Caicos (rows 208-209) and
#   Seychelles/Montserrat (rows 210-211) with refreshed counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 05:36"

# Row 50: now Bolivia (was Barein)
Set-Row 50 @("Bolivia", 19073, 614, 3430, 15011, 0, 21, 632)

# Row 51: now Barein (was Bolivia)
Set-Row 51 @("Barein", 19013, 0, 13267, 5700, 0, 0, 46)

# Row 67: now Honduras (was Marruecos)
Set-Row 67 @("Honduras", 9178, 320, 1025, 7831, 0, 10, 322)

# Row 68: now Marruecos (was Honduras)
Set-Row 68 @("Marruecos", 8885, 0, 7828, 845, 0, 0, 212)

# Row 82: Haiti refreshed counts
Set-Row 82 @("Haiti", 4441, 132, 24, 4341, 0, 3, 76)

# Row 199: Nueva Caledonia refreshed counts
Set-Row 199 @("Nueva Caledonia", 21, 0, 21, 0, 0, 0, 0)

# Row 208: now Santa Sede (was Islas Turcas y Caicos)
Set-Row 208 @("Santa Sede", 12, 0, 12, 0, 0, 0, 0)

# Row 209: now Islas Turcas y Caicos (was Santa Sede)
Set-Row 209 @("Islas Turcas y Caicos", 12, 0, 11, 0, 0, 0, 1)

# Row 210: now Seychelles (was Montserrat)
Set-Row 210 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)

# Row 211: now Montserrat (was Seychelles)
Set-Row 211 @("Montserrat", 11, 0, 10, 0, 0, 0, 1)
